$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 data (row index 9 in the sheet)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 0.63166666666666671
$ws.Range("B9").NumberFormat = "h:mm:ss"

# Set D9 before C9 so the shared-string table gets the "What I did" text
# allocated before the "Comprehension scores" text, matching the source order.
$ws.Range("D9").Value = "Watched Dragon Ball Z, Youtube videos about Minecraft and science. Read Harry Potter and manga."
$ws.Range("C9").Value = "Dragon Ball Z (Audiovisual, Japanese, Familiar):37; Berserk (Text with visuals, Japanese, Familiar):41; ¿Qué Pasó en los Primeros Minutos Después de la Desaparición de los Dinosaurios? [https://www.youtube.com/watch?v=SNSTA4G6O9I] (Audiovisual, English, New):35; Puede que Hayas Nacido Hace 5 Segundos [https://www.youtube.com/watch?v=98kebMYttzg] (Audiovisual, Spanish, New):33;"

# An extra blank but time-formatted cell below, matching the trailing empty row.
$ws.Range("B11").NumberFormat = "[h]:mm:ss"

# Column C got narrower (stored width 106.5703125 characters).
$ws.Columns.Item(3).ColumnWidth = 105.7369791666667
